# Append two new rows (119, 120) of date/remn_amt data to every worksheet
# in the workbook, matching the "Add files via upload" commit.
#
# New rows:
#   row 119 -> date serial 45988 (2025-11-27)
#   row 120 -> date serial 45989 (2025-11-28)
#
# Column B values differ per sheet.

$wb = $excel.ActiveWorkbook

$newRowValues = @{
    1 = @(10608951, 0)
    2 = @(9646317, 0)
    3 = @(3048473, 0)
    4 = @(913355, 0)
    5 = @(1490374, 0)
    6 = @(1775925, 0)
    7 = @(250227, 0)
    8 = @(298827, 0)
}

$dateSerials = @(45988, 45989)

for ($i = 1; $i -le 8; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $values = $newRowValues[$i]

    $ws.Range("A119").Value = $dateSerials[0]
    $ws.Range("B119").Value = $values[0]

    $ws.Range("A120").Value = $dateSerials[1]
    $ws.Range("B120").Value = $values[1]

    # Match the date number format already used in column A (copy down
    # from the last existing data row, A118) on the two new date cells.
    $ws.Range("A119:A120").NumberFormat = $ws.Range("A118").NumberFormat
}
